$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per the crypto price refresh.
# Numeric-looking "Price" values are pinned to Text format before the
# write (then reset to the Normal style) so Excel keeps the exact
# string (e.g. "138.30", "1.00") instead of silently coercing them to
# a number and dropping the trailing zero / losing precision.

$ws.Range("D2").Value = "56.465.42"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "3.004.53"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.431"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.367"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("D12").Value = "3.507.85"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "56.399.51"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "2.997.74"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.497"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("D25").Value = "3.120.79"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.165"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").Value = "0.0₃0941"
$ws.Range("E28").Value = "  +5.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.04%  "
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0663"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "3.032.44"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.653"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").Value = "2.190.50"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.932"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0853"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.65%  "
